# Automatische test-sync: 2025-08-05 18:13:50
# Append a new testmail row to the "Logs" sheet and a matching aggregate
# row to the "Dashboard" sheet, then extend the chart's source ranges.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append row 25 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A25").Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("B25").Value = "mailmind.test@zohomail.eu"
$logs.Range("C25").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("D25").Value = "Inkoop / Bestellingen"
$logs.Range("E25").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F25").Value = "2025-08-05 18:13:03"
$logs.Range("G25").Value = "Ja"
$logs.Range("H25").Value = "Ja"
$logs.Range("I25").Value = "Nee"
$logs.Range("J25").Value = "Nee"

# ---- Dashboard sheet: append row 6 ----------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Inkoop / Bestellingen"
$dash.Range("B6").Value = 1

# ---- Chart: extend category/value source ranges ----------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$series.Values = "='Dashboard'!`$B`$2:`$B`$6"

# ---- Extend conditional-formatting ranges to include the new row 25 -------
$logs.Range("D2:D24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D25"))
$logs.Range("G2:G24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G25"))
$logs.Range("H2:H24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H25"))
$logs.Range("I2:I24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I25"))
$logs.Range("J2:J24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J25"))
